$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.490.85"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "'1.829.53"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -1.17%  "
$ws.Range("D5").Value = "'333.26"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").Value = "'0.3835"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").Value = "'46.16"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").Value = "'0.07865"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "'0.9590"
$ws.Range("E11").Value = "  -3.52%  "
$ws.Range("D12").Value = "'21.07"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").Value = "'1.831.09"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").Value = "'5.840"
$ws.Range("E14").Value = "  -1.39%  "
$ws.Range("D15").Value = "'7.052"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "'89.83"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("D18").Value = "'0.06594"
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("D19").Value = "'0.00001023"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("D20").Value = "'17.12"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "'27.472.02"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "'5.299"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "'2.267"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").Value = "'159.44"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").Value = "'2.040.70"
$ws.Range("E27").Value = "  -2.28%  "
$ws.Range("D28").Value = "'19.38"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("D30").Value = "'5.289"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("D31").Value = "'117.94"
$ws.Range("E31").Value = "  -2.82%  "
$ws.Range("D32").Value = "'0.09381"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").Value = "'0.9305"
$ws.Range("E33").Value = "  -3.92%  "
$ws.Range("D34").Value = "'3.575"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").Value = "'5.223"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "'1.316"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'0.05910"
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("D39").Value = "'8.103"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").Value = "'1.145"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("D42").Value = "'0.5730"
$ws.Range("E42").Value = "  -2.87%  "
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").Value = "'9.905"
$ws.Range("E44").Value = "  -3.56%  "
$ws.Range("D45").Value = "'1.267"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").Value = "'0.5399"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("D48").Value = "'1.888"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "'0.06936"
$ws.Range("E49").Value = "  +3.58%  "
$ws.Range("B50").Value = "PaxosStandard"
$ws.Range("C50").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D50").Value = "'1.046"
$ws.Range("E50").Value = "  -30.02%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'110.43"
$ws.Range("E51").Value = "  -0.74%  "
